$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The tooltip/description text for the "revenu_travail" (Revenu mensuel médian)
# row had a stray double space before "moins." - fix it to a single space.
$ws.Range("D7").Value = "Niveau de salaire ou traitement mensuel net primes incluses médian. Le revenu médian est la valeur telle que la moitié des individus de la population considérée gagne plus, l'autre moitié gagne moins."
